$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Fall")
$ws.Cells.Item(2, 2).Value = 54
$ws.Cells.Item(3, 2).Value = 59
$ws.Cells.Item(4, 2).Value = 64
$ws.Cells.Item(5, 2).Value = 68
$ws.Cells.Item(6, 2).Value = 72
$ws.Cells.Item(7, 2).Value = 75
$ws.Cells.Item(8, 2).Value = 79
$ws.Cells.Item(11, 2).Value = 88
$ws.Cells.Item(12, 2).Value = 91
$ws.Cells.Item(13, 2).Value = 94
$ws.Cells.Item(14, 2).Value = 97
$ws.Cells.Item(15, 2).Value = 100
$ws.Cells.Item(16, 2).Value = 103
$ws.Cells.Item(17, 2).Value = 105
$ws.Cells.Item(18, 2).Value = 108
$ws.Cells.Item(19, 2).Value = 110
$ws.Cells.Item(20, 2).Value = 113
$ws.Cells.Item(21, 2).Value = 115
$ws.Cells.Item(22, 2).Value = 117
$ws.Cells.Item(23, 2).Value = 119
$ws.Cells.Item(24, 2).Value = 122
$ws.Cells.Item(25, 2).Value = 124
$ws.Cells.Item(26, 2).Value = 126
$ws.Cells.Item(27, 2).Value = 128
$ws.Cells.Item(28, 2).Value = 130

$ws = $wb.Worksheets.Item("K-Spring")
$ws.Cells.Item(2, 2).Value = 50
$ws.Cells.Item(3, 2).Value = 54
$ws.Cells.Item(4, 2).Value = 58
$ws.Cells.Item(7, 2).Value = 68
$ws.Cells.Item(8, 2).Value = 72
$ws.Cells.Item(9, 2).Value = 74
$ws.Cells.Item(10, 2).Value = 77
$ws.Cells.Item(11, 2).Value = 80
$ws.Cells.Item(12, 2).Value = 83
$ws.Cells.Item(13, 2).Value = 85
$ws.Cells.Item(14, 2).Value = 87
$ws.Cells.Item(15, 2).Value = 90
$ws.Cells.Item(16, 2).Value = 92
$ws.Cells.Item(17, 2).Value = 94
$ws.Cells.Item(18, 2).Value = 97
$ws.Cells.Item(19, 2).Value = 99
$ws.Cells.Item(20, 2).Value = 101
$ws.Cells.Item(21, 2).Value = 103
$ws.Cells.Item(22, 2).Value = 105
$ws.Cells.Item(23, 2).Value = 107
$ws.Cells.Item(24, 2).Value = 109
$ws.Cells.Item(25, 2).Value = 111
$ws.Cells.Item(27, 2).Value = 114
$ws.Cells.Item(28, 2).Value = 116
$ws.Cells.Item(29, 2).Value = 118
$ws.Cells.Item(30, 2).Value = 120
$ws.Cells.Item(31, 2).Value = 121
$ws.Cells.Item(32, 2).Value = 123
$ws.Cells.Item(33, 2).Value = 125
$ws.Cells.Item(34, 2).Value = 127
$ws.Cells.Item(35, 2).Value = 128

$ws = $wb.Worksheets.Item("1-Fall")
$ws.Cells.Item(2, 2).Value = 44
$ws.Cells.Item(3, 2).Value = 48
$ws.Cells.Item(4, 2).Value = 51
$ws.Cells.Item(5, 2).Value = 54
$ws.Cells.Item(6, 2).Value = 57
$ws.Cells.Item(7, 2).Value = 60
$ws.Cells.Item(8, 2).Value = 63
$ws.Cells.Item(9, 2).Value = 65
$ws.Cells.Item(10, 2).Value = 68
$ws.Cells.Item(11, 2).Value = 70
$ws.Cells.Item(12, 2).Value = 72
$ws.Cells.Item(13, 2).Value = 75
$ws.Cells.Item(14, 2).Value = 77
$ws.Cells.Item(15, 2).Value = 79
$ws.Cells.Item(16, 2).Value = 81
$ws.Cells.Item(17, 2).Value = 83
$ws.Cells.Item(18, 2).Value = 85
$ws.Cells.Item(19, 2).Value = 87
$ws.Cells.Item(20, 2).Value = 89
$ws.Cells.Item(21, 2).Value = 91
$ws.Cells.Item(25, 2).Value = 98
$ws.Cells.Item(26, 2).Value = 100
$ws.Cells.Item(27, 2).Value = 102
$ws.Cells.Item(28, 2).Value = 104
$ws.Cells.Item(29, 2).Value = 106
$ws.Cells.Item(30, 2).Value = 107
$ws.Cells.Item(31, 2).Value = 109
$ws.Cells.Item(32, 2).Value = 111
$ws.Cells.Item(33, 2).Value = 113
$ws.Cells.Item(34, 2).Value = 114
$ws.Cells.Item(35, 2).Value = 116
$ws.Cells.Item(36, 2).Value = 118
$ws.Cells.Item(37, 2).Value = 119
$ws.Cells.Item(38, 2).Value = 121
$ws.Cells.Item(39, 2).Value = 123
$ws.Cells.Item(40, 2).Value = 124

$ws = $wb.Worksheets.Item("1-Spring")
$ws.Cells.Item(2, 2).Value = 40
$ws.Cells.Item(3, 2).Value = 41
$ws.Cells.Item(4, 2).Value = 44
$ws.Cells.Item(5, 2).Value = 47
$ws.Cells.Item(6, 2).Value = 49
$ws.Cells.Item(7, 2).Value = 52
$ws.Cells.Item(8, 2).Value = 54
$ws.Cells.Item(9, 2).Value = 57
$ws.Cells.Item(10, 2).Value = 59
$ws.Cells.Item(11, 2).Value = 61
$ws.Cells.Item(12, 2).Value = 63
$ws.Cells.Item(13, 2).Value = 65
$ws.Cells.Item(14, 2).Value = 67
$ws.Cells.Item(15, 2).Value = 69
$ws.Cells.Item(16, 2).Value = 71
$ws.Cells.Item(17, 2).Value = 73
$ws.Cells.Item(18, 2).Value = 75
$ws.Cells.Item(19, 2).Value = 77
$ws.Cells.Item(20, 2).Value = 79
$ws.Cells.Item(21, 2).Value = 81
$ws.Cells.Item(22, 2).Value = 83
$ws.Cells.Item(23, 2).Value = 85
$ws.Cells.Item(24, 2).Value = 86
$ws.Cells.Item(25, 2).Value = 88
$ws.Cells.Item(26, 2).Value = 90
$ws.Cells.Item(27, 2).Value = 92
$ws.Cells.Item(28, 2).Value = 94
$ws.Cells.Item(29, 2).Value = 96
$ws.Cells.Item(30, 2).Value = 98
$ws.Cells.Item(31, 2).Value = 100
$ws.Cells.Item(32, 2).Value = 102
$ws.Cells.Item(33, 2).Value = 104
$ws.Cells.Item(34, 2).Value = 106
$ws.Cells.Item(35, 2).Value = 108
$ws.Cells.Item(36, 2).Value = 110
$ws.Cells.Item(37, 2).Value = 113
$ws.Cells.Item(38, 2).Value = 115
$ws.Cells.Item(39, 2).Value = 117
$ws.Cells.Item(40, 2).Value = 120

$ws = $wb.Worksheets.Item("2-Fall")
$ws.Cells.Item(6, 2).Value = 42
$ws.Cells.Item(7, 2).Value = 44
$ws.Cells.Item(8, 2).Value = 47
$ws.Cells.Item(9, 2).Value = 49
$ws.Cells.Item(10, 2).Value = 51
$ws.Cells.Item(11, 2).Value = 53
$ws.Cells.Item(12, 2).Value = 54
$ws.Cells.Item(13, 2).Value = 56
$ws.Cells.Item(14, 2).Value = 58
$ws.Cells.Item(15, 2).Value = 60
$ws.Cells.Item(16, 2).Value = 62
$ws.Cells.Item(17, 2).Value = 64
$ws.Cells.Item(18, 2).Value = 65
$ws.Cells.Item(19, 2).Value = 67
$ws.Cells.Item(20, 2).Value = 69
$ws.Cells.Item(21, 2).Value = 71
$ws.Cells.Item(22, 2).Value = 73
$ws.Cells.Item(24, 2).Value = 76
$ws.Cells.Item(25, 2).Value = 78
$ws.Cells.Item(26, 2).Value = 80
$ws.Cells.Item(27, 2).Value = 82
$ws.Cells.Item(28, 2).Value = 84
$ws.Cells.Item(29, 2).Value = 86
$ws.Cells.Item(30, 2).Value = 88
$ws.Cells.Item(31, 2).Value = 90
$ws.Cells.Item(32, 2).Value = 93
$ws.Cells.Item(33, 2).Value = 95
$ws.Cells.Item(34, 2).Value = 98
$ws.Cells.Item(35, 2).Value = 101
$ws.Cells.Item(36, 2).Value = 104
$ws.Cells.Item(38, 2).Value = 113
$ws.Cells.Item(39, 2).Value = 124
$ws.Cells.Item(40, 2).Value = 124

$ws = $wb.Worksheets.Item("2-Spring")
$ws.Cells.Item(9, 2).Value = 41
$ws.Cells.Item(10, 2).Value = 43
$ws.Cells.Item(11, 2).Value = 45
$ws.Cells.Item(12, 2).Value = 46
$ws.Cells.Item(13, 2).Value = 48
$ws.Cells.Item(14, 2).Value = 50
$ws.Cells.Item(15, 2).Value = 51
$ws.Cells.Item(16, 2).Value = 53
$ws.Cells.Item(17, 2).Value = 54
$ws.Cells.Item(18, 2).Value = 56
$ws.Cells.Item(19, 2).Value = 58
$ws.Cells.Item(20, 2).Value = 59
$ws.Cells.Item(21, 2).Value = 61
$ws.Cells.Item(22, 2).Value = 62
$ws.Cells.Item(24, 2).Value = 65
$ws.Cells.Item(25, 2).Value = 67
$ws.Cells.Item(26, 2).Value = 69
$ws.Cells.Item(27, 2).Value = 71
$ws.Cells.Item(28, 2).Value = 72
$ws.Cells.Item(29, 2).Value = 74
$ws.Cells.Item(30, 2).Value = 76
$ws.Cells.Item(31, 2).Value = 78
$ws.Cells.Item(32, 2).Value = 80
$ws.Cells.Item(33, 2).Value = 82
$ws.Cells.Item(34, 2).Value = 85
$ws.Cells.Item(35, 2).Value = 87
$ws.Cells.Item(36, 2).Value = 90
$ws.Cells.Item(37, 2).Value = 93
$ws.Cells.Item(38, 2).Value = 98
$ws.Cells.Item(39, 2).Value = 104
$ws.Cells.Item(40, 2).Value = 114

